# Update "想去人数" (interest count, column F) values across the four
# worksheets of the 杭州-漫展信息 workbook, reflecting newly regenerated
# output data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 27
$ws1.Range("F3").Value = 8254
$ws1.Range("F4").Value = 1933
$ws1.Range("F6").Value = 161
$ws1.Range("F7").Value = 2094
$ws1.Range("F8").Value = 582
$ws1.Range("F9").Value = 54
$ws1.Range("F11").Value = 58
$ws1.Range("F14").Value = 67
$ws1.Range("F15").Value = 7
$ws1.Range("F16").Value = 8625
$ws1.Range("F21").Value = 1817
$ws1.Range("F26").Value = 49
$ws1.Range("F29").Value = 1011
$ws1.Range("F30").Value = 12
$ws1.Range("F32").Value = 13
$ws1.Range("F33").Value = 2133
$ws1.Range("F34").Value = 854
$ws1.Range("F35").Value = 496
$ws1.Range("F39").Value = 210
$ws1.Range("F40").Value = 154
$ws1.Range("F42").Value = 52
$ws1.Range("F44").Value = 49
$ws1.Range("F45").Value = 3975

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 395
$ws2.Range("F3").Value = 208

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2330
$ws3.Range("F3").Value = 712

# --- Sheet 4: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2330
$ws4.Range("F3").Value = 712
$ws4.Range("F4").Value = 395
$ws4.Range("F5").Value = 27
$ws4.Range("F6").Value = 8254
$ws4.Range("F9").Value = 1933
$ws4.Range("F11").Value = 2094
$ws4.Range("F13").Value = 582
$ws4.Range("F17").Value = 58
$ws4.Range("F19").Value = 67
$ws4.Range("F20").Value = 8625
$ws4.Range("F24").Value = 1817
$ws4.Range("F28").Value = 49
$ws4.Range("F31").Value = 12
$ws4.Range("F33").Value = 13
$ws4.Range("F34").Value = 2133
$ws4.Range("F35").Value = 854
$ws4.Range("F37").Value = 496
$ws4.Range("F40").Value = 210
$ws4.Range("F41").Value = 154
$ws4.Range("F44").Value = 3975
